$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 19-94 (B and/or C columns, plus a few E/G columns) per the
# refreshed IsraelStatus dataset pull.
$ws.Range("B19").Value = 183
$ws.Range("C19").Value = 0
$ws.Range("B20").Value = 218
$ws.Range("C20").Value = 0
$ws.Range("B21").Value = 248
$ws.Range("C21").Value = 0
$ws.Range("B22").Value = 251
$ws.Range("C22").Value = 0
$ws.Range("B23").Value = 274
$ws.Range("C23").Value = 0
$ws.Range("B24").Value = 315
$ws.Range("C24").Value = 0
$ws.Range("B25").Value = 361
$ws.Range("C25").Value = 0
$ws.Range("B26").Value = 395
$ws.Range("C26").Value = 0
$ws.Range("B27").Value = 422
$ws.Range("C27").Value = 0
$ws.Range("B28").Value = 450
$ws.Range("C28").Value = 1
$ws.Range("B29").Value = 472
$ws.Range("C29").Value = 1
$ws.Range("B30").Value = 528
$ws.Range("C30").Value = 1
$ws.Range("B31").Value = 606
$ws.Range("C31").Value = 2
$ws.Range("B32").Value = 713
$ws.Range("C32").Value = 2
$ws.Range("B33").Value = 821
$ws.Range("C33").Value = 2
$ws.Range("B34").Value = 992
$ws.Range("C34").Value = 3
$ws.Range("B35").Value = 1110
$ws.Range("C35").Value = 6
$ws.Range("B36").Value = 1270
$ws.Range("C36").Value = 6
$ws.Range("B37").Value = 1460
$ws.Range("C37").Value = 7
$ws.Range("B38").Value = 1631
$ws.Range("C38").Value = 12
$ws.Range("B39").Value = 1865
$ws.Range("C39").Value = 15
$ws.Range("B40").Value = 1994
$ws.Range("C40").Value = 15
$ws.Range("B41").Value = 2167
$ws.Range("C41").Value = 17
$ws.Range("B42").Value = 2442
$ws.Range("C42").Value = 30
$ws.Range("B43").Value = 2780
$ws.Range("C43").Value = 34
$ws.Range("B44").Value = 3202
$ws.Range("C44").Value = 43
$ws.Range("B45").Value = 3721
$ws.Range("C45").Value = 47
$ws.Range("B46").Value = 4242
$ws.Range("C46").Value = 80
$ws.Range("B47").Value = 4645
$ws.Range("C47").Value = 99
$ws.Range("B48").Value = 5252
$ws.Range("C48").Value = 120
$ws.Range("B49").Value = 5948
$ws.Range("C49").Value = 155
$ws.Range("B50").Value = 6561
$ws.Range("C50").Value = 186
$ws.Range("B51").Value = 7745
$ws.Range("C51").Value = 224
$ws.Range("B52").Value = 9090
$ws.Range("C52").Value = 297
$ws.Range("B53").Value = 10747
$ws.Range("C53").Value = 392
$ws.Range("C55").Value = 648
$ws.Range("C56").Value = 818
$ws.Range("C57").Value = 1026
$ws.Range("C58").Value = 1275
$ws.Range("C59").Value = 1623
$ws.Range("C60").Value = 2046
$ws.Range("C61").Value = 2471
$ws.Range("C62").Value = 2998
$ws.Range("B63").Value = 49299
$ws.Range("C63").Value = 3434
$ws.Range("C64").Value = 3916
$ws.Range("C65").Value = 4452
$ws.Range("C66").Value = 4999
$ws.Range("C67").Value = 5732
$ws.Range("C68").Value = 6428
$ws.Range("C69").Value = 7150
$ws.Range("C70").Value = 7748
$ws.Range("C71").Value = 8176
$ws.Range("G71").Value = 49
$ws.Range("C72").Value = 8755
$ws.Range("G72").Value = 55
$ws.Range("B73").Value = 130603
$ws.Range("C73").Value = 9206
$ws.Range("G73").Value = 62
$ws.Range("B74").Value = 137531
$ws.Range("C74").Value = 9585
$ws.Range("G74").Value = 69
$ws.Range("B75").Value = 143401
$ws.Range("C75").Value = 9925
$ws.Range("G75").Value = 78
$ws.Range("B76").Value = 149149
$ws.Range("C76").Value = 10269
$ws.Range("G76").Value = 87
$ws.Range("C77").Value = 10629
$ws.Range("G77").Value = 95
$ws.Range("C78").Value = 10974
$ws.Range("G78").Value = 101
$ws.Range("B79").Value = 173789
$ws.Range("C79").Value = 11534
$ws.Range("G79").Value = 108
$ws.Range("B80").Value = 184748
$ws.Range("C80").Value = 11978
$ws.Range("G80").Value = 116
$ws.Range("B81").Value = 197411
$ws.Range("C81").Value = 12383
$ws.Range("G81").Value = 124
$ws.Range("B82").Value = 206861
$ws.Range("C82").Value = 12694
$ws.Range("G82").Value = 137
$ws.Range("B83").Value = 219645
$ws.Range("C83").Value = 12995
$ws.Range("E83").Value = 182
$ws.Range("G83").Value = 144
$ws.Range("B84").Value = 231206
$ws.Range("C84").Value = 13295
$ws.Range("G84").Value = 153
$ws.Range("B85").Value = 241836
$ws.Range("C85").Value = 13605
$ws.Range("G85").Value = 164
$ws.Range("B86").Value = 253299
$ws.Range("C86").Value = 13896
$ws.Range("G86").Value = 172
$ws.Range("B87").Value = 268646
$ws.Range("C87").Value = 14199
$ws.Range("G87").Value = 178
$ws.Range("B88").Value = 282231
$ws.Range("C88").Value = 14491
$ws.Range("G88").Value = 187
$ws.Range("B89").Value = 295525
$ws.Range("C89").Value = 14720
$ws.Range("G89").Value = 190
$ws.Range("B90").Value = 309683
$ws.Range("C90").Value = 15000
$ws.Range("G90").Value = 194
$ws.Range("B91").Value = 322296
$ws.Range("C91").Value = 15255
$ws.Range("G91").Value = 197
$ws.Range("B92").Value = 331374
$ws.Range("C92").Value = 15415
$ws.Range("G92").Value = 200
$ws.Range("B93").Value = 339841
$ws.Range("C93").Value = 15503
$ws.Range("G93").Value = 202
$ws.Range("B94").Value = 350860
$ws.Range("C94").Value = 15618
$ws.Range("G94").Value = 206

# Append two new rows (95, 96) for 2020-04-20 and 2020-04-21 with the date
# column formatted the same way as the rest of column A.
$ws.Range("A95").Value = 43949
$ws.Range("A95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B95").Value = 361501
$ws.Range("C95").Value = 15786
$ws.Range("D95").Value = 391
$ws.Range("E95").Value = 125
$ws.Range("F95").Value = 98
$ws.Range("G95").Value = 211
$ws.Range("A96").Value = 43950
$ws.Range("A96").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B96").Value = 370505
$ws.Range("C96").Value = 15869
$ws.Range("D96").Value = 370
$ws.Range("E96").Value = 118
$ws.Range("F96").Value = 93
$ws.Range("G96").Value = 217
